$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Verified insulin doses in manuscript and figure legends:
# "mU/kg" -> "U/kg" for the two ITT insulin doses
$d.Content.Find.Execute("0.75 mU/kg (NCD)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "0.75 U/kg (NCD)", 2)
$d.Content.Find.Execute("1.5 mU/kg (HFD)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "1.5 U/kg (HFD)", 2)

# Clamp infusion rate gains a "/min" unit suffix
$d.Content.Find.Execute("/kg following a prime continuous infusion", $false, $false, $false, $false, $false,
                         $true, 1, $false, "/kg/min following a prime continuous infusion", 2)

# Word re-anchors the hidden "_GoBack" bookmark to the location of the last
# edit (right after "Mice we", before "re fasted...").
$full = $d.Content.Text
$idx = $full.IndexOf("Mice we") + 7

$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$r = $d.Range($idx, $idx)
$d.Bookmarks.Add("_GoBack", $r)
